# Adds the new MIGS:virus 4.0 columns (Q..AE) introduced upstream:
# biotic_relationship, encoded_traits, estimated_size, health_state, host,
# host_taxid, isolation_source, num_replicons, pathogenicity, ref_biomaterial,
# samp_collect_device, samp_mat_process, samp_size, source_material_id,
# subspecf_gen_lin -- together with their header-row comments, matching the
# style used by the existing optional (yellow) header cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE")

$headers = @(
  "biotic_relationship",
  "encoded_traits",
  "estimated_size",
  "health_state",
  "host",
  "host_taxid",
  "isolation_source",
  "num_replicons",
  "pathogenicity",
  "ref_biomaterial",
  "samp_collect_device",
  "samp_mat_process",
  "samp_size",
  "source_material_id",
  "subspecf_gen_lin"
)

$comments = @(
  "Free-living or from host (define relationship)",
  "Traits like antibiotic resistance/xenobiotic degration phenotypes/converting phage genes",
  "Estimated size of genome",
  "Health or disease status of sample at time of collection",
  "The natural (as opposed to laboratory) host to the organism from which the sample was obtained. Use the full taxonomic name, eg, `"Homo sapiens`".",
  "NCBI taxonomy ID of the host, e.g. 9606",
  "Describes the physical, environmental and/or local geographical source of the biological sample from which the sample was derived.",
  "Reports the number of replicons in a nuclear genome of eukaryotes, in the genome of a bacterium or archaea or the number of segments in a segmented virus. Always applied to the haploid chromosome count of a eukaryote",
  "To what is the entity pathogenic",
  "Primary publication or genome report in the form of pubmed ID, DOI or URL",
  "Method or device employed for collecting sample",
  "Processing applied to the sample during or after isolation",
  "Amount or size of sample (volume, mass or area) that was collected",
  "unique identifier assigned to a material sample used for extracting nucleic acids, and subsequent sequencing. The identifier can refer either to the original material collected or to any derived sub-samples.",
  "Information about the genetic distinctness of the lineage (eg., biovar, serovar)"
)

# Copy the yellow "optional field" style from the existing C15 header cell
# (sample_title's sibling "description") onto the whole new range in one go,
# then fill in the text and per-cell comments.
$ws.Range("C15").Copy()
$ws.Range("Q15:AE15").PasteSpecial(-4122)

for ($i = 0; $i -lt $cols.Length; $i++) {
  $cell = $ws.Range($cols[$i] + "15")
  $cell.Value = $headers[$i]
  $cell.AddComment($comments[$i])
}
